$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) values re-pulled from source data / mean calculation update.
# Map of row -> new value for column F
$updates = @{
    2  = -2
    3  = -4
    4  = -3
    5  = 1
    6  = -3
    7  = 3
    8  = -3
    9  = 13
    10 = 3
    11 = 1
    12 = -4
    13 = -2
    14 = -4
    15 = 2
    16 = -3
    18 = 4
    20 = -2
    21 = 2
    22 = 3
    24 = 5
    25 = 1
    26 = 2
    28 = 2
    29 = -2
    30 = -2
    31 = 1
    32 = -1
    33 = 6
    34 = -1
    37 = -1
    38 = 2
    40 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
